# StudyHistory.xlsx update
# - Drop the "ID" column (E) entirely
# - Replace the sample rows with the new test data (4 data rows now)
# - Dates are stored as plain DD/MM/YYYY text, not Excel serials
# - Clear the old datetime number-format from column D so the new text
#   values are not affected by the previous style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover datetime number-format from column D before we
# overwrite the values, so the new text dates carry no special format.
$ws.Range("D2:D3").ClearFormats()

# Drop the "ID" column (E1:E4) - everything right of it shifts left and
# the used-range dimension shrinks from E to D automatically.
$ws.Columns("E").Delete()

# Row 2
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = "test subject"
$ws.Range("C2").Value = "test summary"
$ws.Range("D2").Value = "16/06/2025"

# Row 3
$ws.Range("A3").Value = 4
$ws.Range("D3").Value = "16/06/2025"

# Row 4
$ws.Range("A4").Value = 25
$ws.Range("D4").Value = "16/06/2025"

# Row 5 (new row)
$ws.Range("A5").Value = 5
$ws.Range("D5").Value = "16/06/2025"

# Subject/Summary are blanked out on rows 3-5, but keep the cells present
# (touching a formatting property - reset to its own default - stops the
# engine from dropping the now-empty cell from the sheet entirely).
foreach ($addr in @("B3", "C3", "B4", "C4", "B5", "C5")) {
    $ws.Range($addr).Value = ""
    $ws.Range($addr).Font.Bold = $false
}

Write-Output "StudyHistory sheet updated"
